# Added KNN k=7 results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the Test Accuracy (column D) values for the KNN k=7 block (rows 42-49)
$ws.Range("D42").Value = 0.72729999999999995
$ws.Range("D43").Value = 0.77270000000000005
$ws.Range("D44").Value = 0.77270000000000005
$ws.Range("D45").Value = 0.72729999999999995
$ws.Range("D46").Value = 0.72729999999999995
$ws.Range("D47").Value = 0.81820000000000004
$ws.Range("D48").Value = 0.81820000000000004
$ws.Range("D49").Value = 0.86360000000000003

# Match the styling used by the other "Test Accuracy" column cells
# (General number format, centered, 12pt font - same as the rest of column D)
$range = $ws.Range("D42:D49")
$range.HorizontalAlignment = -4108
$range.Font.Size = 12

# Update the active selection to reflect where the user last clicked after entering data
$ws.Range("D50").Select()
